# RPA datasets push 2024-05-31
# 1) Rename the 시프트업 (ShiftUp) listing to indicate its new KOSPI (유가) listing.
# 2) Fill in the previously blank "확정공모가" (confirmed offering price) for
#    디비금융스팩12호 (row 20) with 2000 -- stored as text, matching the
#    existing "-" placeholders used in the rest of that column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "시프트업(유가)"

# Force the new value into column D (확정공모가) as text "2000", same as the
# other cells in that column (all stored as text, e.g. "-"), then drop back
# to the default "Normal" style so no stray per-cell number format sticks.
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2000"
$ws.Range("D20").Style = "Normal"
